$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 159 (shifts existing rows 159-185 down to 162-188),
# mirroring a new week of "Terminal La Palmera de La Serena - Chirimoya" data being
# prepended to the weekly log.
$ws.Range("A159:A161").EntireRow.Insert()

# New row 159: Especial, Provincia de Limarí, $/bandeja 10 kilos
$ws.Cells.Item(159,1).Value = 8
$ws.Cells.Item(159,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(159,3).Value = "Coquimbo"
$ws.Cells.Item(159,4).Value = [DateTime]::FromOADate(44841)
$ws.Cells.Item(159,5).Value = 4
$ws.Cells.Item(159,6).Value = "Fruta"
$ws.Cells.Item(159,7).Value = 100107
$ws.Cells.Item(159,8).Value = "Otros"
$ws.Cells.Item(159,9).Value = 100107002
$ws.Cells.Item(159,10).Value = "Chirimoya"
$ws.Cells.Item(159,11).Value = "Cultivar IV Región"
$ws.Cells.Item(159,12).Value = "Especial"
$ws.Cells.Item(159,13).Value = 360
$ws.Cells.Item(159,14).Value = 22000
$ws.Cells.Item(159,15).Value = 23000
$ws.Cells.Item(159,16).Value = 22500
$ws.Cells.Item(159,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(159,18).Value = "Provincia de Limarí"
$ws.Cells.Item(159,19).Value = 2250
$ws.Cells.Item(159,20).Value = 10

# New row 160: Primera, Provincia de Limarí, $/bandeja 10 kilos
$ws.Cells.Item(160,1).Value = 8
$ws.Cells.Item(160,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(160,3).Value = "Coquimbo"
$ws.Cells.Item(160,4).Value = [DateTime]::FromOADate(44841)
$ws.Cells.Item(160,5).Value = 4
$ws.Cells.Item(160,6).Value = "Fruta"
$ws.Cells.Item(160,7).Value = 100107
$ws.Cells.Item(160,8).Value = "Otros"
$ws.Cells.Item(160,9).Value = 100107002
$ws.Cells.Item(160,10).Value = "Chirimoya"
$ws.Cells.Item(160,11).Value = "Cultivar IV Región"
$ws.Cells.Item(160,12).Value = "Primera"
$ws.Cells.Item(160,13).Value = 400
$ws.Cells.Item(160,14).Value = 19000
$ws.Cells.Item(160,15).Value = 20000
$ws.Cells.Item(160,16).Value = 19500
$ws.Cells.Item(160,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(160,18).Value = "Provincia de Limarí"
$ws.Cells.Item(160,19).Value = 1950
$ws.Cells.Item(160,20).Value = 10

# New row 161: Segunda, Provincia de Limarí, $/bandeja 10 kilos
$ws.Cells.Item(161,1).Value = 8
$ws.Cells.Item(161,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(161,3).Value = "Coquimbo"
$ws.Cells.Item(161,4).Value = [DateTime]::FromOADate(44841)
$ws.Cells.Item(161,5).Value = 4
$ws.Cells.Item(161,6).Value = "Fruta"
$ws.Cells.Item(161,7).Value = 100107
$ws.Cells.Item(161,8).Value = "Otros"
$ws.Cells.Item(161,9).Value = 100107002
$ws.Cells.Item(161,10).Value = "Chirimoya"
$ws.Cells.Item(161,11).Value = "Cultivar IV Región"
$ws.Cells.Item(161,12).Value = "Segunda"
$ws.Cells.Item(161,13).Value = 400
$ws.Cells.Item(161,14).Value = 15000
$ws.Cells.Item(161,15).Value = 16000
$ws.Cells.Item(161,16).Value = 15500
$ws.Cells.Item(161,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(161,18).Value = "Provincia de Limarí"
$ws.Cells.Item(161,19).Value = 1550
$ws.Cells.Item(161,20).Value = 10
